$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.601.64"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.15"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.12"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.812.14"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.588.58"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.76"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.631.38"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.64"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("E23").Value = "  -4.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.36"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.47"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.651"
$ws.Range("E34").Value = "  +22.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.304.05"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0172"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.817"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.786"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.28"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.92"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.724.39"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.86"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.831"
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.00%  "
